$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 originally held the sheet-name placeholder string; remove it so the
# header text can occupy C1 and the shared-string table is rebuilt cleanly.
$ws.Range("A1").ClearContents()

$ws.Range("C1").Value = "%(n)of students who think exacerbation`nN = 1179"

$ws.Range("A2").Value = "idep.diff.know.students"
$ws.Range("B2").Value = "Difficulty getting to know other students in class "
$ws.Range("C2").Value = "61.74% (728)"

$ws.Range("A3").Value = "idep.diff.help.students"
$ws.Range("B3").Value = "Difficulty getting help from other students in class "
$ws.Range("C3").Value = "49.53% (584)"

$ws.Range("A4").Value = "idep.talk.unknown"
$ws.Range("B4").Value = "Needing to talk with students who I don’t know during online group work"
$ws.Range("C4").Value = "43.60% (514)"

$ws.Range("A5").Value = "idep.diff.know.instructor"
$ws.Range("B5").Value = "Difficulty getting to know instructors"
$ws.Range("C5").Value = "51.31% (605)"

$ws.Range("A6").Value = "idep.diff.help.instructor"
$ws.Range("B6").Value = "Difficulty getting help from instructors"
$ws.Range("C6").Value = "47.41% (559)"

$ws.Range("A7").Value = "idep.compare"
$ws.Range("B7").Value = "Comparing myself to other students"
$ws.Range("C7").Value = "45.89% (541)"

$ws.Range("A8").Value = "idep.not.inperson"
$ws.Range("B8").Value = "Not having to show up in person to online science courses "
$ws.Range("C8").Value = "30.28% (357)"

$ws.Range("A9").Value = "idep.pace"
$ws.Range("B9").Value = "Deciding the pace at which I work through an online science course"
$ws.Range("C9").Value = "27.48% (324)"

$ws.Range("A10").Value = "idep.camera"
$ws.Range("B10").Value = " Being on camera"
$ws.Range("C10").Value = "39.27% (463)"

$ws.Range("A11").Value = "idep.nav.tech"
$ws.Range("B11").Value = "Needing to navigate technology in high-pressure situations (e.g., during exams)"
$ws.Range("C11").Value = "45.29% (534"

$ws.Range("A12").Value = "idep.proctor.exam"
$ws.Range("B12").Value = "Online monitored proctored testing"
$ws.Range("C12").Value = "57.51% (678)"

$ws.Range("A13").Value = "idep.questions"
$ws.Range("B13").Value = "Struggling to have questions answered "
$ws.Range("C13").Value = "48.43% (571)"

$ws.Range("A14").Value = "idep.comm.instructor"
$ws.Range("B14").Value = "Struggling to communicate effectively with the instructor"
$ws.Range("C14").Value = "46.73% (551)"

$ws.Range("A15").Value = "idep.home.distract"
$ws.Range("B15").Value = "At-home distractions that can interfere with online science courses"
$ws.Range("C15").Value = "54.37% (641)"

$ws.Range("A16").Value = "idep.personal.tech"
$ws.Range("B16").Value = "The potential for personal technology issues (e.g., unstable internet connection) "
$ws.Range("C16").Value = "43.94% (518)"

$ws.Range("A17").Value = "idep.other"
$ws.Range("A18").Value = "idep.nothing"

# Column B (the long survey-item descriptions) uses Times New Roman 12pt
$ws.Range("B2:B16").Font.Name = "Times New Roman"
$ws.Range("B2:B16").Font.Size = 12

# Wrap text across the whole used range
$ws.Range("A2:C18").WrapText = $true
$ws.Range("B1:C1").WrapText = $true

# Row heights
$ws.Rows(1).RowHeight = 40.5
$ws.Rows("2:16").RowHeight = 15.75

# Column widths (characters). The runtime stores widths using an internal
# 7-px max-digit-width model with a fixed 5px padding, so the input value
# is pre-adjusted (target - 5/7) to land as close as possible to the
# OOXML width recorded in the reference file.
$ws.Columns("A").ColumnWidth = 27.285714285714285
$ws.Columns("B").ColumnWidth = 59.660714285714285
$ws.Columns("C").ColumnWidth = 22.660714285714285

$ws.Range("D10").Select()